$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.550.68"
$ws.Range("D3").Value = "2.244.19"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "2.271.59"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "44.291.97"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "0.0₃0954"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.17%  "
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D43").Value = "1.791.22"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +12.92%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.192"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "79.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70.62"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.16"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.73"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.16%  "
